$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.456404
$ws.Range("H2").Value = 1.369212
$ws.Range("I2").Value = 0.01914960767004715
$ws.Range("J2").Value = 0.01914960767004715
$ws.Range("M2").Value = 0.6537256666666666
$ws.Range("N2").Value = 1.961177
$ws.Range("O2").Value = 0.08287237534104652
$ws.Range("P2").Value = 0.08287237534104651
$ws.Range("Q2").Value = 0.2983630091693333
$ws.Range("R2").Value = 2.685267082524
$ws.Range("S2").Value = 0.001586973474465931
$ws.Range("T2").Value = 0.001586973474465931
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.456404
$ws.Range("H3").Value = 1.369212
$ws.Range("I3").Value = 0.01914960767004715
$ws.Range("J3").Value = 0.01914960767004715
$ws.Range("O3").Value = 0.7131728182689164
$ws.Range("P3").Value = 0.7131728182689163
$ws.Range("Q3").Value = 2.567615411538667
$ws.Range("R3").Value = 23.108538703848
$ws.Range("S3").Value = 0.01365697967079159
$ws.Range("T3").Value = 0.01365697967079158
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.456404
$ws.Range("H4").Value = 1.369212
$ws.Range("I4").Value = 0.01914960767004715
$ws.Range("J4").Value = 0.01914960767004715
$ws.Range("O4").Value = 0.2039548063900371
$ws.Range("P4").Value = 0.2039548063900371
$ws.Range("Q4").Value = 0.7342925735946667
$ws.Range("R4").Value = 6.608633162352
$ws.Range("S4").Value = 0.003905654524789636
$ws.Range("T4").Value = 0.003905654524789636
$ws.Range("I5").Value = 0.8285024587002443
$ws.Range("J5").Value = 0.8285024587002443
$ws.Range("M5").Value = 0.6537256666666666
$ws.Range("N5").Value = 1.961177
$ws.Range("O5").Value = 0.08287237534104652
$ws.Range("P5").Value = 0.08287237534104651
$ws.Range("Q5").Value = 12.90859274723656
$ws.Range("R5").Value = 116.177334725129
$ws.Range("S5").Value = 0.06865996672838653
$ws.Range("T5").Value = 0.06865996672838653
$ws.Range("I6").Value = 0.8285024587002443
$ws.Range("J6").Value = 0.8285024587002443
$ws.Range("O6").Value = 0.7131728182689164
$ws.Range("P6").Value = 0.7131728182689163
$ws.Range("S6").Value = 0.5908654334139798
$ws.Range("T6").Value = 0.5908654334139797
$ws.Range("I7").Value = 0.8285024587002443
$ws.Range("J7").Value = 0.8285024587002443
$ws.Range("O7").Value = 0.2039548063900371
$ws.Range("P7").Value = 0.2039548063900371
$ws.Range("S7").Value = 0.168977058557878
$ws.Range("T7").Value = 0.168977058557878
$ws.Range("I8").Value = 0.1523479336297086
$ws.Range("J8").Value = 0.1523479336297086
$ws.Range("M8").Value = 0.6537256666666666
$ws.Range("N8").Value = 1.961177
$ws.Range("O8").Value = 0.08287237534104652
$ws.Range("P8").Value = 0.08287237534104651
$ws.Range("Q8").Value = 2.373677241941
$ws.Range("R8").Value = 21.363095177469
$ws.Range("S8").Value = 0.01262543513819406
$ws.Range("T8").Value = 0.01262543513819406
$ws.Range("I9").Value = 0.1523479336297086
$ws.Range("J9").Value = 0.1523479336297086
$ws.Range("O9").Value = 0.7131728182689164
$ws.Range("P9").Value = 0.7131728182689163
$ws.Range("S9").Value = 0.1086504051841451
$ws.Range("T9").Value = 0.1086504051841451
$ws.Range("I10").Value = 0.1523479336297086
$ws.Range("J10").Value = 0.1523479336297086
$ws.Range("O10").Value = 0.2039548063900371
$ws.Range("P10").Value = 0.2039548063900371
$ws.Range("S10").Value = 0.03107209330736945
$ws.Range("T10").Value = 0.03107209330736944
